$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Consolidate the "This library aims..." paragraph's many runs
#    (split around the smart quotes / CSS / HTML words) into a
#    single run containing the same text. Re-assigning identical
#    text via .Text is treated as a no-op, so delete the range and
#    retype it to force the run to be rebuilt as one piece.
# ------------------------------------------------------------------
$para = $d.Paragraphs(3)
$paraRange = $para.Range
$bodyRange = $d.Range($paraRange.Start, $paraRange.End - 1)

$fullText = "This library aims to enable users to easily display content in the form of " + `
    [char]0x201C + "tiles" + [char]0x201D + `
    ", which are customizable rectangular blocks which can hold content of different types, such as images, text, links, etc. " + `
    "Developers have the option of sorting and customizing the tiles, including things like color, size, position, opacity, etc. " + `
    "End users will be able to interact with the content in the tiles, depending on what developers have put there. " + `
    "There could also be things like customizable buttons which can be mapped to different functionality. " + `
    "Developers would use this instead of writing their own code because its often frustrating to get the CSS and HTML code just right to align different blocks in the way that you intend to. " + `
    "I hope that with this library it cuts down on the frustration of dealing with the structure of the page and allow developers to focus their time on producing the actual content that would go on there."

$bodyRange.Delete()
$insPoint = $d.Range($bodyRange.Start, $bodyRange.Start)
$insPoint.InsertAfter($fullText)

# ------------------------------------------------------------------
# 2) Append a new run of text right after the hyperlink in the
#    "My alpha release is deployed at:" paragraph, describing the
#    alpha release's current functionality.
# ------------------------------------------------------------------
$linkRange = $d.Content
$linkRange.Find.Execute("https://csc309-alpha-release.herokuapp.com/")
$linkRange.Collapse(0)
$linkRange.InsertAfter(". The page just has a basic use case of the library, which uses the library to create some tiles and some tiles with images. The tiles show a visual indication when they are being hovered over, and the developer can specify what they want the tile to do when it is clicked.")
